$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows (2-10) with the new "BF-C2DL-HSC" primary-
#     track parameter set, and append new data rows (11-14) for the newly
#     added datasets. Column A also gets re-alphabetized/expanded.
#
# Columns B-E and G-J hold numeric-looking text (stored as shared strings in
# the original file, not real numbers), so those cells are explicitly
# formatted as Text before the value is written, to stop Excel from auto-
# converting the text into a floating point number. Column F is a genuine
# number (1) and column A holds plain dataset-name text that Excel already
# stores as text without help.

# Row 2: BF-C2DL-HSC
$ws.Range("A2").Value = "BF-C2DL-HSC"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "23.9718"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "4.1848"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.93934"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.0072582"
$ws.Range("F2").Value = 1
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0.026126"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "8.4398"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "189.8257"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "5.6477"

# Row 3: BF-C2DL-MuSC
$ws.Range("A3").Value = "BF-C2DL-MuSC"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "23.9718"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "4.1848"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0.93934"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.0072582"
$ws.Range("F3").Value = 1
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "0.026126"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "8.4398"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "189.8257"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "5.6477"

# Row 4: DIC-C2DH-HeLa
$ws.Range("A4").Value = "DIC-C2DH-HeLa"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "23.9718"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "4.1848"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.93934"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.0072582"
$ws.Range("F4").Value = 1
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "0.026126"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "8.4398"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "189.8257"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "5.6477"

# Row 5: Fluo-C2DL-MSC
$ws.Range("A5").Value = "Fluo-C2DL-MSC"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "23.9718"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "4.1848"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.93934"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.0072582"
$ws.Range("F5").Value = 1
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "0.026126"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "8.4398"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "189.8257"
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = "5.6477"

# Row 6: Fluo-C3DH-A549
$ws.Range("A6").Value = "Fluo-C3DH-A549"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "23.9718"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "4.1848"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.93934"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.0072582"
$ws.Range("F6").Value = 1
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "0.026126"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "8.4398"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "189.8257"
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "5.6477"

# Row 7: Fluo-C3DH-H157
$ws.Range("A7").Value = "Fluo-C3DH-H157"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "23.9718"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "4.1848"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.93934"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.0072582"
$ws.Range("F7").Value = 1
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "0.026126"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "8.4398"
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "189.8257"
$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value = "5.6477"

# Row 8: Fluo-C3DL-MDA231
$ws.Range("A8").Value = "Fluo-C3DL-MDA231"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "23.9718"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "4.1848"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.93934"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.0072582"
$ws.Range("F8").Value = 1
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "0.026126"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "8.4398"
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "189.8257"
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "5.6477"

# Row 9: Fluo-N2DH-GOWT1
$ws.Range("A9").Value = "Fluo-N2DH-GOWT1"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "23.9718"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "4.1848"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.93934"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.0072582"
$ws.Range("F9").Value = 1
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "0.026126"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "8.4398"
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "189.8257"
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = "5.6477"

# Row 10: Fluo-N2DL-HeLa
$ws.Range("A10").Value = "Fluo-N2DL-HeLa"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "23.9718"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "4.1848"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.93934"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.0072582"
$ws.Range("F10").Value = 1
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "0.026126"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "8.4398"
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "189.8257"
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J10").Value = "5.6477"

# Row 11: Fluo-N3DH-CE
$ws.Range("A11").Value = "Fluo-N3DH-CE"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "23.9718"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "4.1848"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.93934"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.0072582"
$ws.Range("F11").Value = 1
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "0.026126"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "8.4398"
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "189.8257"
$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value = "5.6477"

# Row 12: Fluo-N3DH-CHO
$ws.Range("A12").Value = "Fluo-N3DH-CHO"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "23.9718"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "4.1848"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.93934"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.0072582"
$ws.Range("F12").Value = 1
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "0.026126"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "8.4398"
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "189.8257"
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = "5.6477"

# Row 13: PhC-C2DH-U373
$ws.Range("A13").Value = "PhC-C2DH-U373"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "23.9718"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "4.1848"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.93934"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.0072582"
$ws.Range("F13").Value = 1
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "0.026126"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "8.4398"
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "189.8257"
$ws.Range("J13").NumberFormat = "@"
$ws.Range("J13").Value = "5.6477"

# Row 14: PhC-C2DL-PSC
$ws.Range("A14").Value = "PhC-C2DL-PSC"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "23.9718"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "4.1848"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.93934"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.0072582"
$ws.Range("F14").Value = 1
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0.026126"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "8.4398"
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "189.8257"
$ws.Range("J14").NumberFormat = "@"
$ws.Range("J14").Value = "5.6477"

# --- Clear the temporary Text number-format so the cells fall back to the
#     default (General) style, matching the target formatting exactly ---
$ws.Range("B2:E14").Style = "Normal"
$ws.Range("G2:J14").Style = "Normal"

# --- Match the saved selection/active cell from the target file ---
$ws.Range("A3").Select()
